$d = $word.ActiveDocument

# 1) Merge the "Spēlētājs izvēlas ... doto skaitli:" runs into one sentence.
#    (Text content is unchanged, just re-wording across runs collapses to
#    the same visible sentence, so a direct Find/Replace on the full phrase
#    is sufficient and safe.)
$d.Content.Find.Execute(
    "Spēlētājs izvēlas grūtības pakāpi ierakstos pakāpes doto skaitli:",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "Spēlētājs izvēlas grūtības pakāpi ierakstos pakāpes doto skaitli:",
    2)

# 2) Extend the "Katrs jautājums tiek nummurēts;" sentence.
$d.Content.Find.Execute(
    "Katrs jautājums tiek nummurēts;",
    $false, $true, $false, $false, $false,
    $true, 1, $false,
    "Katrs jautājums tiek nummurēts un atrodas apākšā iepriekšējam jautājumam. Jautājumi veido kolonnu;",
    2)

# 3) Insert a new list paragraph after the "...(Skatīt tabulā pašā apakšā);" item.
$targetIndex = 0
$i = 0
foreach ($p in $d.Paragraphs) {
    $i = $i + 1
    if ($p.Range.Text -like "*Skatīt tabulā pašā apakšā*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$target.Range.InsertParagraphAfter()

$newPara = $d.Paragraphs.Item($targetIndex + 1)
$newPara.Range.Text = "Viens pareizi atbildēts jautājums sniedz 1 punktu, bet punkti netiek zaudēti, ja ir atbildēts nepareizi;"
